$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B17 was stored as a text "3" -- correct it to a genuine number.
$ws.Cells.Item(17, 2).Value = 3

# Append a new annotation row (row 18) for Sunsi Wu.
$ws.Cells.Item(18, 1).Value = "Sunsi Wu"
# politeness_score "4" is entered as text (leading apostrophe keeps it a string).
$ws.Cells.Item(18, 2).Value = "'4"
$ws.Cells.Item(18, 3).Value = "thank"
$ws.Cells.Item(18, 4).Value = "ACK"
$ws.Cells.Item(18, 5).Value = "OTH"
$ws.Cells.Item(18, 6).Value = "e27c53be-a9c3-4697-b8f8-f90bcc73c090"
$ws.Cells.Item(18, 7).Value = "SJaP_-xAb_annotated.xlsx"
$ws.Cells.Item(18, 8).Value = "We thank everybody again for their useful suggestions and we uploaded a revision of the paper."
